$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the color labels in column D, aligned with rows 15-19 and 22
$ws.Range("D15").Value = "rot"
$ws.Range("D16").Value = "blau"
$ws.Range("D17").Value = "braun"
$ws.Range("D18").Value = "gelb"
$ws.Range("D19").Value = "orange"
$ws.Range("D22").Value = "grün"

# Update the Max values in column C for the two panels that were re-wired
$ws.Range("C18").Value = 52
$ws.Range("C19").Value = 53
$ws.Range("C22").Value = 51

# Sync the view position / selection to match where the user was working
$ws.Range("C18").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
